$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: merge the two runs
#   "This rule is basically the opposite of the last one,"
#   " it takes a stop as an input, and outputs the number of lines"
# into a single run, while leaving the following run
#   " that passes through it."
# untouched (it must stay a separate <w:r>).
#
# A plain Find/Replace over the combined text would normally also
# absorb the following same-formatted run into the merge, so we
# temporarily flip a character formatting flag on that trailing run
# to break the run-coalescing chain, perform the merge, then restore
# the original formatting.
# -----------------------------------------------------------------

$tail = $d.Content
$tail.Find.Execute(" that passes through it.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Font.Bold = 1

$merge = $d.Content
$merge.Find.Execute( `
    "This rule is basically the opposite of the last one, it takes a stop as an input, and outputs the number of lines", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "This rule is basically the opposite of the last one, it takes a stop as an input, and outputs the number of lines", `
    2)

$tail2 = $d.Content
$tail2.Find.Execute(" that passes through it.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail2.Font.Bold = 0

# -----------------------------------------------------------------
# Change 2: after the run
#   "D – Anyway, thank you very much for listening, do you have any questions?"
# add a new run containing a single space, right before the
# _GoBack bookmark.
#
# Inserting directly at the end of the existing run would simply
# extend that run's text instead of creating a new <w:r>, so the new
# space is given a momentary formatting nudge (bold on/off) to force
# it into its own run, matching the target structure.
# -----------------------------------------------------------------

$q = $d.Content
$q.Find.Execute("D – Anyway, thank you very much for listening, do you have any questions?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$qEnd = $q.End

$ins = $d.Range($qEnd, $qEnd)
$ins.InsertAfter(" ")

$newRun = $d.Range($qEnd, $qEnd + 1)
$newRun.Font.Bold = 1
$newRun.Font.Bold = 0
